# zqlexamples.xlsx - "fixing axis labels other minor changes"
#
# The sheet has several small 4-5 column "mini tables" that each start
# with a header row of f1 / x1<-{...} / y1<-{...} / z1<-'state'.* labels
# (see rows 37 and 41). The bottom of the sheet had two free-floating
# commentary notes ("^ works perfectly" in D40 and the long "this query
# does not work..." note in E44) instead of a proper labelled row for the
# last mini-table - this adds the missing axis-label header row (45) and
# removes the stray notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two leftover commentary cells - they're being replaced by a
# real header row below.
$ws.Range("D40").Value = ""
$ws.Range("E44").Clear()

# Row 45: new axis-label header row for the last mini-table, matching the
# same look (style/borders/fill) as the other f1/x1/y1/z1 header rows
# (row 37 / row 41).
$ws.Range("A37:D37").Copy($ws.Range("A45:D45"))
$ws.Range("A45").Value = "f1"
$ws.Range("B45").Value = "x1<-{'year','month'}"
$ws.Range("C45").Value = "y1<-{'soldprice','listingprice'}"
$ws.Range("D45").Value = "z1<-'state'.*"

# Minor row-height touch-ups that came with the edit.
$ws.Rows.Item(39).RowHeight = 24
$ws.Rows.Item(44).RowHeight = 16
$ws.Rows.Item(45).RowHeight = 16

# Leave the selection on the newly-added header row.
$ws.Range("A45:D45").Select()
